$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 46265.086
$ws.Range("I28").Value = 56605.89
$ws.Range("J28").Value = 9038.200000000001
$ws.Range("K28").Value = 56605.89
$ws.Range("L28").Value = 9038.200000000001
$ws.Range("M28").Value = -56120.89
$ws.Range("N28").Value = -10008.2
$ws.Range("H32").Value = 3180.9443
$ws.Range("I32").Value = 1538.7142
$ws.Range("J32").Value = 4226
$ws.Range("K32").Value = 1538.7142
$ws.Range("L32").Value = 4226
$ws.Range("M32").Value = -1212.7142
$ws.Range("N32").Value = -4878
$ws.Range("H40").Value = 3147.2632
$ws.Range("I40").Value = 2916.8333
$ws.Range("J40").Value = 3253.6155
$ws.Range("K40").Value = 2916.8333
$ws.Range("L40").Value = 3253.6155
$ws.Range("M40").Value = -2741.8333
$ws.Range("N40").Value = -3603.6155
$ws.Range("H51").Value = 7869.231
$ws.Range("I51").Value = 5350.25
$ws.Range("J51").Value = 8988.777
$ws.Range("K51").Value = 5350.25
$ws.Range("L51").Value = 8988.777
$ws.Range("M51").Value = -4866.25
$ws.Range("N51").Value = -9956.777
$ws.Range("H52").Value = 636.7143
$ws.Range("I52").Value = 87.40000000000001
$ws.Range("K52").Value = 262.2
$ws.Range("M52").Value = -102.2
$ws.Range("H106").Value = 3657.1853
$ws.Range("I106").Value = 1701.4348
$ws.Range("K106").Value = 1701.4348
$ws.Range("M106").Value = -1070.4348
$ws.Range("H113").Value = 166670080
$ws.Range("J113").Value = 4101.2
$ws.Range("L113").Value = 4101.2
$ws.Range("N113").Value = -10609.2
$ws.Range("H132").Value = 2360.3052
$ws.Range("I132").Value = 2038.3948
$ws.Range("K132").Value = 6115.1844
$ws.Range("M132").Value = -3585.1844
$ws.Range("H137").Value = 16669136
$ws.Range("I137").Value = 62502224
$ws.Range("J137").Value = 2558.0454
$ws.Range("K137").Value = 187506672
$ws.Range("L137").Value = 7674.1362
$ws.Range("M137").Value = -187504122
$ws.Range("N137").Value = -12774.1362
$ws.Range("H138").Value = 3107.8928
$ws.Range("I138").Value = 2703.9
$ws.Range("J138").Value = 3332.3333
$ws.Range("K138").Value = 8111.700000000001
$ws.Range("L138").Value = 9996.999899999999
$ws.Range("M138").Value = -2971.700000000001
$ws.Range("N138").Value = -20276.9999
$ws.Range("H141").Value = 3013.7334
$ws.Range("I141").Value = 3025.5833
$ws.Range("J141").Value = 2966.3333
$ws.Range("K141").Value = 9076.749899999999
$ws.Range("L141").Value = 8898.999899999999
$ws.Range("M141").Value = -3896.749899999999
$ws.Range("N141").Value = -19258.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1512.7084
$ws.Range("I2").Value = 1372.7646
$ws.Range("K2").Value = 1372.7646
$ws.Range("M2").Value = -1259.7646
$ws.Range("H32").Value = 345.77
$ws.Range("I32").Value = 326.6979
$ws.Range("J32").Value = 803.5
$ws.Range("K32").Value = 326.6979
$ws.Range("L32").Value = 803.5
$ws.Range("M32").Value = -39.6979
$ws.Range("N32").Value = -1377.5
$ws.Range("H116").Value = 1512.7084
$ws.Range("I116").Value = 1372.7646
$ws.Range("K116").Value = 1372.7646
$ws.Range("M116").Value = 921.2354

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1512.7084
$ws.Range("I3").Value = 1372.7646
$ws.Range("K3").Value = 1372.7646
$ws.Range("M3").Value = -1258.7646
$ws.Range("H56").Value = 29998
$ws.Range("J56").Value = 29998
$ws.Range("L56").Value = 29998
$ws.Range("N56").Value = -31476
$ws.Range("H134").Value = 1683.8334
$ws.Range("I134").Value = 1544
$ws.Range("K134").Value = 4632
$ws.Range("M134").Value = -2097

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 466.64706
$ws.Range("I7").Value = 388.07144
$ws.Range("K7").Value = 388.07144
$ws.Range("M7").Value = -275.07144
$ws.Range("H16").Value = 2488.6
$ws.Range("I16").Value = 2209.5557
$ws.Range("K16").Value = 2209.5557
$ws.Range("M16").Value = -1922.5557
$ws.Range("H19").Value = 610.7143
$ws.Range("I19").Value = 670.8333
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 670.8333
$ws.Range("L19").Value = 250
$ws.Range("M19").Value = -500.8333
$ws.Range("N19").Value = -590
$ws.Range("H24").Value = 610.7143
$ws.Range("I24").Value = 670.8333
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 670.8333
$ws.Range("L24").Value = 250
$ws.Range("M24").Value = -500.8333
$ws.Range("N24").Value = -590
$ws.Range("H31").Value = 2804.074
$ws.Range("J31").Value = 3126.9688
$ws.Range("L31").Value = 3126.9688
$ws.Range("N31").Value = -3716.9688
$ws.Range("H34").Value = 2804.074
$ws.Range("J34").Value = 3126.9688
$ws.Range("L34").Value = 3126.9688
$ws.Range("N34").Value = -3530.9688
$ws.Range("H51").Value = 69249.836
$ws.Range("I51").Value = 68999.75
$ws.Range("J51").Value = 69750
$ws.Range("K51").Value = 68999.75
$ws.Range("L51").Value = 69750
$ws.Range("M51").Value = -68263.75
$ws.Range("N51").Value = -71222
$ws.Range("H60").Value = 15000
$ws.Range("I60").Value = 15000
$ws.Range("K60").Value = 15000
$ws.Range("M60").Value = -14489
$ws.Range("H61").Value = 69249.836
$ws.Range("I61").Value = 68999.75
$ws.Range("J61").Value = 69750
$ws.Range("K61").Value = 68999.75
$ws.Range("L61").Value = 69750
$ws.Range("M61").Value = -68651.75
$ws.Range("N61").Value = -70446
$ws.Range("H99").Value = 2871.3076
$ws.Range("I99").Value = 2482.8333
$ws.Range("J99").Value = 3204.2856
$ws.Range("K99").Value = 2482.8333
$ws.Range("L99").Value = 3204.2856
$ws.Range("M99").Value = -984.8332999999998
$ws.Range("N99").Value = -6200.2856
$ws.Range("H113").Value = 2488.6
$ws.Range("I113").Value = 2209.5557
$ws.Range("K113").Value = 2209.5557
$ws.Range("M113").Value = -39.55569999999989
$ws.Range("H126").Value = 2871.3076
$ws.Range("I126").Value = 2482.8333
$ws.Range("J126").Value = 3204.2856
$ws.Range("K126").Value = 7448.499899999999
$ws.Range("L126").Value = 9612.856800000001
$ws.Range("M126").Value = -4978.499899999999
$ws.Range("N126").Value = -14552.8568
$ws.Range("H132").Value = 2668.5
$ws.Range("I132").Value = 1803.421
$ws.Range("K132").Value = 5410.263
$ws.Range("M132").Value = -2880.263
$ws.Range("H134").Value = 2416
$ws.Range("I134").Value = 1715.15
$ws.Range("K134").Value = 5145.450000000001
$ws.Range("M134").Value = -2610.450000000001
$ws.Range("H141").Value = 58749.332
$ws.Range("I141").Value = 39999
$ws.Range("J141").Value = 62499.4
$ws.Range("K141").Value = 39999
$ws.Range("L141").Value = 62499.4
$ws.Range("M141").Value = -34819
$ws.Range("N141").Value = -72859.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1103.3334
$ws.Range("I108").Value = 1103.3334
$ws.Range("K108").Value = 3310.0002
$ws.Range("M108").Value = -430.0001999999999
$ws.Range("H122").Value = 624.75
$ws.Range("J122").Value = 750
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1001550
$ws.Range("I11").Value = 1001550
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1001550
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1001411
$ws.Range("N11").ClearContents()
$ws.Range("H132").Value = 404137.16
$ws.Range("I132").Value = 529896.3
$ws.Range("J132").Value = 5899.8335
$ws.Range("K132").Value = 1589688.9
$ws.Range("L132").Value = 17699.5005
$ws.Range("M132").Value = -1587158.9
$ws.Range("N132").Value = -22759.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27781152
$ws.Range("I7").Value = 45457160
$ws.Range("J7").Value = 4571.2856
$ws.Range("K7").Value = 45457160
$ws.Range("L7").Value = 4571.2856
$ws.Range("M7").Value = -45457048
$ws.Range("N7").Value = -4795.2856
$ws.Range("H25").Value = 30007
$ws.Range("I25").Value = 30007
$ws.Range("K25").Value = 30007
$ws.Range("M25").Value = -29777
$ws.Range("H40").Value = 3771.1428
$ws.Range("I40").Value = 3088.5557
$ws.Range("K40").Value = 3088.5557
$ws.Range("M40").Value = -2952.5557
$ws.Range("H46").Value = 2344.9429
$ws.Range("I46").Value = 1732.2222
$ws.Range("J46").Value = 2993.7058
$ws.Range("K46").Value = 1732.2222
$ws.Range("L46").Value = 2993.7058
$ws.Range("M46").Value = -1544.2222
$ws.Range("N46").Value = -3369.7058
$ws.Range("H51").Value = 25000
$ws.Range("I51").Value = 25000
$ws.Range("K51").Value = 25000
$ws.Range("M51").Value = -24522
$ws.Range("H126").Value = 27781152
$ws.Range("I126").Value = 45457160
$ws.Range("J126").Value = 4571.2856
$ws.Range("K126").Value = 136371480
$ws.Range("L126").Value = 13713.8568
$ws.Range("M126").Value = -136369010
$ws.Range("N126").Value = -18653.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J3").Value = 15000000
$ws.Range("L3").Value = 15000000
$ws.Range("N3").Value = -15000228
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H81").Value = 4516.1665
$ws.Range("J81").Value = 5952.909
$ws.Range("L81").Value = 11905.818
$ws.Range("N81").Value = -14027.818
$ws.Range("H84").Value = 4516.1665
$ws.Range("J84").Value = 5952.909
$ws.Range("L84").Value = 59529.09
$ws.Range("N84").Value = -70137.09
$ws.Range("H126").Value = 2788.6365
$ws.Range("I126").Value = 2899.5
$ws.Range("J126").Value = 2655.6
$ws.Range("K126").Value = 8698.5
$ws.Range("L126").Value = 7966.799999999999
$ws.Range("M126").Value = -6228.5
$ws.Range("N126").Value = -12906.8
$ws.Range("H136").Value = 3989.75
$ws.Range("I136").Value = 2988.2273
$ws.Range("K136").Value = 8964.6819
$ws.Range("M136").Value = -6414.6819
